$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "29.535.04"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  -2.53%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.001.32"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  -4.14%  "

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.014"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  +1.17%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "329.31"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -3.95%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.012"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +0.99%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.5003"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -4.39%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4216"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  -4.52%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "54.24"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  -0.54%  "

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09008"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  -3.38%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.117"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  -4.29%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "23.26"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -6.15%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.047.21"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  -5.46%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.024"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  -6.70%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.464"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  -6.20%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.012"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +0.91%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "94.41"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  -6.67%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00001114"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  -3.86%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06677"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +0.30%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.66"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -6.90%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.011"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +0.79%  "

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.966"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  -5.76%  "

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "29.582.01"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -2.52%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.98"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -4.33%  "

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.301"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  -0.11%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "158.96"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  -2.16%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "20.71"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -4.88%  "

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.351"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  -4.75%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.295"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  -8.65%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "128.15"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -3.66%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.056"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -6.84%  "

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09962"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -4.60%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.564"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -6.14%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.831"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -6.40%  "

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.797"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  -1.59%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.02465"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  -6.10%  "

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.268"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  -8.73%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06415"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  -6.21%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.304"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -3.22%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6534"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -6.25%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.67"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -6.66%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2045"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -7.36%  "

$ws.Range("E43").Value = "  +0.87%  "

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6356"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -6.72%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.52"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  -6.06%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.195"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  -5.87%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.306"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -4.96%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.514"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -3.36%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00000000333"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -4.41%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06992"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  -3.20%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.127"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -6.71%  "
